# Update the "南宁-漫展信息" workbook with refreshed 想去人数/最低票价 figures.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 16287
$ws1.Range("G2").Value = 62
$ws1.Range("F3").Value = 347
$ws1.Range("F4").Value = 723
$ws1.Range("F5").Value = 251
$ws1.Range("F6").Value = 676
$ws1.Range("F7").Value = 1672
$ws1.Range("F8").Value = 155

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 12

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 16287
$ws4.Range("G2").Value = 62
$ws4.Range("F3").Value = 347
$ws4.Range("F4").Value = 723
$ws4.Range("F5").Value = 251
$ws4.Range("F6").Value = 12
$ws4.Range("F8").Value = 676
$ws4.Range("F9").Value = 1672
$ws4.Range("F11").Value = 155
